$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the default (Normal) style so we can restore it after forcing
# a text number-format on cells whose new value would otherwise be
# re-interpreted by Excel as a number (losing formatting like trailing zeros).
$normalStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '60.964.95'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '2.882.14'
$ws.Range("E3").Value = '  -1.54%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.36'
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.38'
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = '  -5.65%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.491'
$ws.Range("D8").Style = $normalStyle
$ws.Range("E8").Value = '  -3.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.87'
$ws.Range("D9").Style = $normalStyle
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("E10").Value = '  -5.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.427'
$ws.Range("D11").Style = $normalStyle
$ws.Range("E11").Value = '  -3.29%  '
$ws.Range("E12").Value = '  -4.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.17'
$ws.Range("D13").Style = $normalStyle
$ws.Range("E13").Value = '  -4.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.126'
$ws.Range("D14").Style = $normalStyle
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").Value = '3.362.62'
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("D16").Value = '60.905.11'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '2.880.73'
$ws.Range("E17").Value = '  -1.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.47'
$ws.Range("D18").Style = $normalStyle
$ws.Range("E18").Value = '  -3.54%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '424.58'
$ws.Range("D19").Style = $normalStyle
$ws.Range("E19").Value = '  -1.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.21'
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = '  -1.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.651'
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = '  -4.17%  '
$ws.Range("E22").Value = '  -2.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.71'
$ws.Range("D23").Style = $normalStyle
$ws.Range("E23").Value = '  -2.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.32'
$ws.Range("D24").Style = $normalStyle
$ws.Range("E24").Value = '  -5.52%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -7.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.34'
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = '  -4.68%  '
$ws.Range("E28").Value = '  -3.62%  '
$ws.Range("E29").Value = '  -9.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.58'
$ws.Range("D30").Style = $normalStyle
$ws.Range("E30").Value = '  -6.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = '  +0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.54'
$ws.Range("D32").Style = $normalStyle
$ws.Range("E32").Value = '  -4.20%  '
$ws.Range("E33").Value = '  -5.49%  '
$ws.Range("D34").Value = '0.0₃0839'
$ws.Range("E34").Value = '  -2.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.967'
$ws.Range("D35").Style = $normalStyle
$ws.Range("E35").Value = '  -4.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.40'
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = '  -4.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.78'
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = '  -7.61%  '
$ws.Range("E39").Value = '  -5.18%  '
$ws.Range("E40").Value = '  -3.09%  '
$ws.Range("E41").Value = '  -6.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.263'
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = '  -7.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '37.98'
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = '  -6.82%  '
$ws.Range("D44").Value = '2.662.56'
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '131.72'
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = '  -1.70%  '
$ws.Range("E46").Value = '  -4.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '346.93'
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = '  -8.93%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.101'
$ws.Range("D49").Style = $normalStyle
$ws.Range("E49").Value = '  -4.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.18'
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = '  -6.93%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.91'
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = '  -4.98%  '
